$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '28.768.53'
$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  +2.68%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.876.22'
$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  +2.55%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '1.005'
$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  +0.38%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '325.67'
$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  +0.41%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '1.005'
$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  +0.36%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.4595'
$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  -0.87%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.3872'
$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  +0.34%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.07856'
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.9910'
$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  +3.39%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '21.81'
$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  -0.27%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '1.875.75'
$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  +1.92%  '
$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  +1.81%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '5.719'
$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  +0.86%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.06949'
$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  +1.40%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '88.52'
$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  +0.29%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '1.006'
$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  +0.41%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '0.00001005'
$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  +1.42%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '16.81'
$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  +0.98%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '1.004'
$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  +0.30%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '28.783.11'
$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  +2.71%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '11.05'
$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  +0.88%  '
$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  +1.02%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.122.93'
$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  +3.14%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '153.16'
$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  -0.83%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '19.23'
$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  +0.56%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '5.823'
$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  +2.74%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '1.973'
$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  +0.74%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '119.17'
$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  +0.62%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '0.09322'
$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  +0.87%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.9191'
$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  -1.60%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '5.310'
$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  +1.24%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '1.341'
$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  +1.93%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '3.326'
$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  +0.63%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.05770'
$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  -0.78%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '1.152'
$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  +1.64%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.02078'
$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  -1.87%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '7.687'
$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  -0.07%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.5651'
$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  +1.26%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.1791'
$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  +1.96%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '9.895'
$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  +0.35%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.07214'
$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  -1.70%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '11.81'
$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  +1.95%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.5300'
$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  +0.93%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.165'
$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  +3.97%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '1.125'
$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  -0.57%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '113.79'
$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  +1.02%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '1.829'
$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  +0.35%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '2.410'
$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  +3.82%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '1.005'
$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  +0.39%  '
